$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------
# The sheet lists variables considered for the weighting model. This
# edit:
#   - removes the "education_degree" row from its original spot
#   - re-adds it (unchanged) further down, together with two brand new
#     variables ("assessment_center", "ethnic_background") that are
#     also excluded from the prediction ("no" / "NA").
#
# We build the three new rows *before* deleting the original
# education_degree row (row 5) so that we can still borrow its
# formatting (Calibri Light font) for the re-added row. Deleting row 5
# then shifts everything below up by one, landing the new rows exactly
# on 17-19 and leaving the pre-existing stray formatted cell on row 20
# (previously row 21), matching the target layout.
# --------------------------------------------------------------------

# New style for the "assessment_center" ID cell: same Menlo font as
# used elsewhere in the sheet (copied from E4), but recoloured to
# FFCE9178 and right aligned.
$ws.Range("E4").Copy()
$ws.Range("A18").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A18").Font.Color = 7901646   # RGB(0xCE,0x91,0x78) packed as BGR
$ws.Range("A18").HorizontalAlignment = -4152 # xlRight

# Row 18 (pre-delete) -> assessment_center
$ws.Range("A18").Value = 54
$ws.Range("B18").Value = "assessment_center"
$ws.Range("C18").Value = "UK Biobank assessment centre"
$ws.Range("D18").Value = "no"
$ws.Range("E18").Value = "NA"
$ws.Range("E4").Copy()
$ws.Range("E18").PasteSpecial(-4122) # xlPasteFormats

# Row 19 (pre-delete) -> ethnic_background
$ws.Range("A2").Copy()
$ws.Range("A19").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A19").Value = 21000
$ws.Range("B19").Value = "ethnic_background"
$ws.Range("C19").Value = "Ethnic background"
$ws.Range("D19").Value = "no"
$ws.Range("E19").Value = "NA"
$ws.Range("E4").Copy()
$ws.Range("E19").PasteSpecial(-4122) # xlPasteFormats

# Row 20 (pre-delete, currently a lone styled stray cell) -> re-added
# education_degree row, borrowing the Calibri Light look of the
# original row 5 before it disappears.
$ws.Range("B1").Copy()
$ws.Range("A20").PasteSpecial(-4122) # xlPasteFormats (plain/default look)
$ws.Range("A20").Value = 6138
$ws.Range("B5").Copy()
$ws.Range("B20").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B20").Value = "education_degree"
$ws.Range("C20").Value = "Education (degree)"
$ws.Range("D20").Value = "no"
$ws.Range("E20").Value = "NA"
$ws.Range("E4").Copy()
$ws.Range("E20").PasteSpecial(-4122) # xlPasteFormats

# Remove the original "education_degree" row - everything below moves
# up by one, putting the three rows built above on 17, 18 and 19, and
# the untouched stray formatted cell back on row 20.
$ws.Rows(5).Delete()

$ws.Range("E17").Select()
